$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 186, shifting the existing rows (186-189) down to (187-190)
$ws.Rows("186").Insert()

# Populate the newly inserted row 186 with the new weekly price record
$ws.Range("A186").Value = 3
$ws.Range("B186").Value = "Femacal de La Calera"
$ws.Range("C186").Value = "Coquimbo"
$ws.Range("D186").Value = 44448
$ws.Range("E186").Value = 5
$ws.Range("F186").Value = 100112043
$ws.Range("G186").Value = "Pepino ensalada"
$ws.Range("H186").Value = "Sin especificar"
$ws.Range("I186").Value = "Primera"
$ws.Range("J186").Value = 150
$ws.Range("K186").Value = 15000
$ws.Range("L186").Value = 15500
$ws.Range("M186").Value = 15233
$ws.Range("N186").Value = "$/caja 70 unidades"
$ws.Range("O186").Value = "Región de Arica y Parinacota"
$ws.Range("P186").Value = 218
$ws.Range("Q186").Value = 70
$ws.Range("R186").Value = "Hortaliza"
